$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column G (St_group): every data row becomes "IP" (rows 3 & 6 were "BOP").
# ---------------------------------------------------------------------------
$ws.Range("G2:G7").Value = "IP"

# ---------------------------------------------------------------------------
# Column B (SLNO running numbers) is removed entirely for the data rows.
# ---------------------------------------------------------------------------
$ws.Range("B2:B7").ClearContents()

# ---------------------------------------------------------------------------
# Column H (St_category "PS") is cleared. Rows 2-4 keep their existing cell
# style (border/number format) but become blank; rows 5-7 had no style so the
# cells disappear completely once cleared.
# ---------------------------------------------------------------------------
$ws.Range("H2:H7").ClearContents()

# ---------------------------------------------------------------------------
# Column M (Rate): rows 2 and 5 change from 100 to 10.
# ---------------------------------------------------------------------------
$ws.Range("M2").Value = 10
$ws.Range("M5").Value = 10

# ---------------------------------------------------------------------------
# Column N (Amount) becomes a computed formula (Quantity * Rate) instead of a
# hard-coded number. N2 is entered on its own (standalone formula) and then
# N3:N7 are entered together so the engine groups them into the same "shared
# formula" group that Excel itself would create when the formula from N3 is
# filled down to N7.
# ---------------------------------------------------------------------------
$ws.Range("N2").Formula = "=J2*M2"
$ws.Range("N3:N7").Formula = "=J3*M3"

# ---------------------------------------------------------------------------
# Column J got a lot wider (used to just fit "Quantity").
# ---------------------------------------------------------------------------
$ws.Columns(10).ColumnWidth = 12.3

# ---------------------------------------------------------------------------
# View state: the sheet is now scrolled so column C is left-most, and the
# selection highlights the St_group column that was just edited.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("G2:G7").Select()
$excel.ActiveWindow.ScrollColumn = 3
